$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.724.05'
$ws.Range('E2').Value = '  +3.98%  '
$ws.Range('D3').Value = '1.914.54'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.84%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.12'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5192'
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3973'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08516'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.93'
$ws.Range('E10').Value = '  +2.81%  '
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.304'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '1.908.71'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.04'
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001117'
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.97'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.035'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').Value = '29.724.11'
$ws.Range('E23').Value = '  +3.80%  '
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.211'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = '2.126.73'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.02'
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.35'
$ws.Range('E28').Value = '  -1.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.470'
$ws.Range('E29').Value = '  +4.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.87'
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.090'
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.200'
$ws.Range('E33').Value = '  +6.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.693'
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02503'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06645'
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.157'
$ws.Range('E37').Value = '  +2.80%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2212'
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.248'
$ws.Range('E39').Value = '  +4.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.209'
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6553'
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.241'
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.41'
$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6140'
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.24'
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.700'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.069'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.239'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.88'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.190'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.52'
$ws.Range('E51').Value = '  +1.90%  '

Write-Host "Applied cryptos update"
